$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 2 days,
# keeping the intraday fraction (and therefore formatting/style) unchanged.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 2
}

# Retrained model output: updated "Actual Production (MW)" values for the
# Horeco solar ramp-up window (rows 30-43).
$newProduction = @{
    30 = 0
    31 = 4
    32 = 13
    33 = 28
    34 = 41
    35 = 56
    36 = 73
    37 = 90
    38 = 104
    39 = 126
    40 = 149
    41 = 163
    42 = 176
    43 = 190
}

foreach ($r in $newProduction.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $newProduction[$r]
}
